$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add contribution details for Marcos Fernandez's rows (10-13)
$ws.Range("C10").Value = "Created dashboard to see sales and profits given certain regions"
$ws.Range("C11").Value = "Added screenshots and description of app functionality with the help of all team members"
$ws.Range("C12").Value = "Description of team reflections on what we learned during this project and the most interesting parts"
$ws.Range("C13").Value = "Final formatting details and cohesion check"

# Update the active selection to C12 as last edited cell
$ws.Range("C12").Select() | Out-Null
